$wb = $excel.ActiveWorkbook

# --- 1. Add the new "cfop hours" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "cfop hours"

# --- 2. Add the new "cfop" column (G) to the "PI hours" sheet ---
$ws1 = $wb.Worksheets.Item("PI hours")

$ws1.Range("G1").Value = "cfop"
$ws1.Range("G2").Value = "['cfop_HUTCHINSON']"
$ws1.Range("G3").Value = "['cfop_PARK']"
$ws1.Range("G4").Value = "['cfop_MITRA']"

# copy the header formatting (bold, border, centered) from an existing header cell
$ws1.Range("F1").Copy()
$ws1.Range("G1").PasteSpecial(-4122)

# --- 3. Populate the new "cfop hours" sheet, mirroring the layout of the
#        other per-group sheets (e.g. "department hours") ---
$ws4.Range("B1").Value = "cfop"
$ws4.Range("C1").Value = "hours"
$ws4.Range("D1").Value = "percentage"

$ws4.Range("A2").Value = 0
$ws4.Range("B2").Value = "cfop_HUTCHINSON"
$ws4.Range("C2").Value = 226
$ws4.Range("D2").Value = 79.5774647887324

$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = "cfop_PARK"
$ws4.Range("C3").Value = 44
$ws4.Range("D3").Value = 15.49295774647887

$ws4.Range("A4").Value = 2
$ws4.Range("B4").Value = "cfop_MITRA"
$ws4.Range("C4").Value = 14
$ws4.Range("D4").Value = 4.929577464788732

# copy header + index-column formatting from the "PI hours" sheet
$ws1.Range("B1:D1").Copy()
$ws4.Range("B1:D1").PasteSpecial(-4122)

$ws1.Range("A2:A4").Copy()
$ws4.Range("A2:A4").PasteSpecial(-4122)

# keep "PI hours" as the active/selected sheet, matching the original workbook
$ws1.Activate()
